$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (appears on Overview!E2:F3 and on the zh-cn/de-de Status column C2:C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Hyperlink colour used by the existing "Latest Target/Handback File" links
# (matches the pre-existing HyperLink cell style: underline + RGB(0x64,0x95,0xED))
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: fill in "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" now that the handback has happened.
# ---------------------------------------------------------------------------
$wsZhCn.Range("I2").Value = "ad31a5f0-6dea-496a-835e-962f7c08b9ac.md"
$wsZhCn.Range("J2").Value = "ad31a5f0-6dea-496a-835e-962f7c08b9ac.6f0eea6a380f5e2b4075a3721895492eabf55f1a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-30 18:31:45"

$wsZhCn.Range("I3").Value = "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md"
$wsZhCn.Range("J3").Value = "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.91c6c5264385abd71ed4a12392f2d8180ae100c1.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-30 18:31:45"

# Recreate the hyperlinks in document order (A2, I2, A3, I3) so relationship
# ids come out the same way the handback tool lays them down, then restore
# the plain (non-theme) hyperlink font used throughout this workbook.
$zhCnAd31Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5a5765a226ccc4caa81076602c639c925bcf7d2/e2e/ad31a5f0-6dea-496a-835e-962f7c08b9ac.md"
$zhCnC1553Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5a5765a226ccc4caa81076602c639c925bcf7d2/e2e/c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnAd31Url, "", "", "ad31a5f0-6dea-496a-835e-962f7c08b9ac.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhCnAd31Url, "", "", "ad31a5f0-6dea-496a-835e-962f7c08b9ac.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhCnC1553Url, "", "", "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhCnC1553Url, "", "", "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md") | Out-Null

$wsZhCn.Range("A2").Font.Underline = 2
$wsZhCn.Range("A2").Font.Color = $hyperlinkColor
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor
$wsZhCn.Range("A3").Font.Underline = 2
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# 3) de-de sheet: same shape of update, handed back a little later.
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = "ad31a5f0-6dea-496a-835e-962f7c08b9ac.md"
$wsDeDe.Range("J2").Value = "ad31a5f0-6dea-496a-835e-962f7c08b9ac.6f0eea6a380f5e2b4075a3721895492eabf55f1a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-30 18:31:52"

$wsDeDe.Range("I3").Value = "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md"
$wsDeDe.Range("J3").Value = "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.91c6c5264385abd71ed4a12392f2d8180ae100c1.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-30 18:31:52"

$deDeAd31Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5a5765a226ccc4caa81076602c639c925bcf7d2/e2e/ad31a5f0-6dea-496a-835e-962f7c08b9ac.md"
$deDeC1553Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5a5765a226ccc4caa81076602c639c925bcf7d2/e2e/c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeAd31Url, "", "", "ad31a5f0-6dea-496a-835e-962f7c08b9ac.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deDeAd31Url, "", "", "ad31a5f0-6dea-496a-835e-962f7c08b9ac.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deDeC1553Url, "", "", "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deDeC1553Url, "", "", "c1553d09-8e1a-4bb1-8b59-3ade70b6c0d2.md") | Out-Null

$wsDeDe.Range("A2").Font.Underline = 2
$wsDeDe.Range("A2").Font.Color = $hyperlinkColor
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor
$wsDeDe.Range("A3").Font.Underline = 2
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# 4) Column widths, widened to fit the new longer status text / filenames.
# ---------------------------------------------------------------------------
$wsOverview.Range("E:E").ColumnWidth = 29.9777047293527
$wsOverview.Range("F:F").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C:C").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I:I").ColumnWidth = 40
$wsZhCn.Range("J:J").ColumnWidth = 40

$wsDeDe.Range("C:C").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I:I").ColumnWidth = 40
$wsDeDe.Range("J:J").ColumnWidth = 40
